$d = $word.ActiveDocument

# The document currently contains a single empty paragraph that holds only
# the "_GoBack" bookmark. The target edit inserts the text "joel" at the
# very start of that paragraph (before the bookmark), with both the new
# run and the paragraph mark itself tagged as French (France) language.

# Step 1: Mark the (still empty) insertion point / paragraph mark as
# French. Doing this before any text is inserted produces the
# <w:pPr><w:rPr><w:lang w:val="fr-FR"/></w:rPr></w:pPr> on the paragraph.
$word.Selection.LanguageID = "fr-FR"

# Step 2: Insert "joel" at the start of the (first) paragraph, ahead of
# the existing bookmark start/end.
$para = $d.Paragraphs(1)
$r = $para.Range
$r.Collapse(1)          # wdCollapseStart
$r.InsertBefore("joel")

# Step 3: Tag the newly-inserted run's text as French as well.
$runRange = $d.Range(0, 4)
$runRange.LanguageID = "fr-FR"
